# Updated Maven TestNG Integration
# - Sheet1: remove the "Vish / Password@123 / Vishwanath / Vishwanath D B" row (row 3)
# - credentials: remove the "Admin / OYs6MbnC2@ / Aaron / Aaron Update" row (old row 3),
#   which shifts the "Vish" row up to row 3, then update its password/name columns
#   to the new test data (free@123 / Fun / Joy) and re-point its hyperlink.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Drop the hyperlink anchored on B3 before the row disappears.
$ws1.Hyperlinks.Delete()

# Remove the whole "Vish" row; rows below (none here) shift up.
$ws1.Rows("3:3").Delete()

$ws1.Range("I21").Select()

# ---------------------------------------------------------------------------
# credentials
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("credentials")
$ws2.Activate()

# Drop only the hyperlink that currently sits on B4 (the "Vish" row); the one
# on B2 ("Tabby" row) is untouched.
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$4') {
        $h.Delete()
    }
}

# Remove the "Admin / OYs6MbnC2@ / Aaron / Aaron Update" row; the "Vish" row
# (old row 4) shifts up to become row 3.
$ws2.Rows("3:3").Delete()

# Refresh the (shifted) "Vish" row with the new credential data.
$ws2.Range("B3").Value = "free@123"
$ws2.Range("C3").Value = "Fun"
$ws2.Range("D3").Value = "Joy"

# That row no longer carries the old "plain" override font - reset to default.
$ws2.Range("A3").Style = "Normal"
$ws2.Range("C3").Style = "Normal"
$ws2.Range("D3").Style = "Normal"

# Re-create the hyperlink on B3 against the new address.
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:free@123")

$ws2.Range("D3").Select()
